$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B (ME) through G (N)
$data = @{
    2  = @(-0.1653664247801546, 1.120513518117024, 9.67144629664304, 3.109894901221429, 3.139811074079776, 46)
    3  = @(0.01719828832093468, 1.063375132475499, 6.205273630094869, 2.491038664913668, 2.519126847305234, 45)
    4  = @(0.005100588602229281, 0.9255207531993044, 4.773114308188849, 2.184745822330106, 2.209997815720953, 44)
    5  = @(0.0498243245536589, 0.9581495051347919, 5.307203266713742, 2.303736805000463, 2.330455680691624, 43)
    6  = @(0.0885671167821773, 0.9573066994189174, 5.181437973596104, 2.276277218090122, 2.302124899960172, 42)
    7  = @(0.1065348192791066, 0.929811024553358, 4.72585632215702, 2.173903475814191, 2.198265070640741, 41)
    8  = @(0.1420065762956791, 0.9765351640082575, 4.879971543114854, 2.209065762514745, 2.232580586974871, 40)
    9  = @(0.1340186904858539, 0.9571447338883464, 4.907955779113894, 2.21539066060907, 2.240240800137693, 39)
    10 = @(0.1600587551659434, 0.9793142779403741, 5.040124732968332, 2.245022212132506, 2.269368439699161, 38)
    11 = @(0.1574861501701002, 0.9407490158730698, 5.047038324643313, 2.246561444662334, 2.271946963475703, 37)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
